# Handles float input without breaking stuff
# ------------------------------------------------------------------
# This marksheet previously showed an "Absent" placeholder because the
# student's answers (and the derived scoring numbers) had not been
# filled in, and the per-question marking penalty ("-1") had been
# stored as plain text instead of a real number, which broke the
# score computation whenever it was combined with numeric inputs.
# This script fills in the student's actual answers, recomputes the
# summary numbers, fixes the penalty cell to be numeric, and removes
# the now-unused third answer block (columns G/H) together with the
# rows of the second answer block (columns D/E) beyond question 3.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Summary block (rows 10-12) ---------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 13
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 52
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "48/112"

# --- Remove the now-unused third answer block (columns G & H) --------
$ws.Range("G15:H40").Clear()

# --- Section 1 answers (column A), rows 16-40 -------------------------
$section1 = @{
    16 = @("correctStyle",   "Option A")
    17 = @("correctStyle",   "Option D")
    18 = @("correctStyle",   "Option B")
    19 = @("correctStyle",   "Option C")
    20 = @("normalStyle",    $null)
    21 = @("incorrectStyle", "Option B")
    22 = @("normalStyle",    $null)
    23 = @("normalStyle",    $null)
    24 = @("normalStyle",    $null)
    25 = @("correctStyle",   "Option A")
    26 = @("correctStyle",   "Option C")
    27 = @("incorrectStyle", "Option C")
    28 = @("incorrectStyle", "Option B")
    29 = @("correctStyle",   "Option D")
    30 = @("correctStyle",   "Option B")
    31 = @("normalStyle",    $null)
    32 = @("normalStyle",    $null)
    33 = @("correctStyle",   "Option D")
    34 = @("normalStyle",    $null)
    35 = @("correctStyle",   "Option D")
    36 = @("normalStyle",    $null)
    37 = @("normalStyle",    $null)
    38 = @("normalStyle",    $null)
    39 = @("correctStyle",   "Option D")
    40 = @("incorrectStyle", "Option B")
}

foreach ($row in $section1.Keys) {
    $styleName = $section1[$row][0]
    $text = $section1[$row][1]
    $cell = $ws.Range("A$row")
    if ($text) {
        $cell.Value = $text
    }
    $cell.Style = $styleName
}

# --- Section 2 answers (columns D/E), only questions 1-3 remain -------
$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"

$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"

# Rows 19-40 of the second answer block (columns D/E) are no longer used
$ws.Range("D19:E40").Clear()
